# Update "想去人数" (want-to-go count) values in the F column
# for sheets "展览" and "全部类型", as generated by the gh-pages build
# at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 73
$ws1.Range("F6").Value = 106
$ws1.Range("F7").Value = 340
$ws1.Range("F8").Value = 4285
$ws1.Range("F10").Value = 4915
$ws1.Range("F11").Value = 548
$ws1.Range("F12").Value = 1234

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 73
$ws4.Range("F6").Value = 106
$ws4.Range("F8").Value = 340
$ws4.Range("F9").Value = 4285
$ws4.Range("F11").Value = 4915
$ws4.Range("F12").Value = 548
$ws4.Range("F13").Value = 1234
